# Production tech level-up tables: add a new "INT_buildTime" (build time,
# in seconds) column in column G for every tech sheet, finishing the dev
# of the production techs.

$wb = $excel.ActiveWorkbook

# Build-time values (column G, rows 2-16) - identical progression on
# every sheet.
$buildTimes = @(34128, 34128, 45504, 56880, 68256, 91008, 113760, 136512, 159264, 182016, 204768, 204768, 204768, 204768, 204768)

$sheetNames = @("crane", "fastFix", "reinforcing", "stoneCarving", "ironSmelting", "seniorTower", "forestation", "cropResearch", "beerSupply")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Extend the used range out to column I (matches H, which already
    # carries the style but no value) so the sheet dimension grows from
    # H16 to I16, same as Excel does when a column gets touched.
    $ws.Range("H1").Copy($ws.Range("I1"))

    # New header label for the build-time column.
    $ws.Range("G1").Value = "INT_buildTime"

    # Fill in the per-level build times.
    for ($i = 0; $i -lt $buildTimes.Length; $i++) {
        $row = $i + 2
        $ws.Range("G$row").Value = $buildTimes[$i]
    }

    # Leave the cursor on G2, like in the saved file.
    $ws.Range("G2").Select()
}
